# Edit script implementing the diff: adds w:proofErr markers (gramStart/
# gramEnd and spellStart/spellEnd) around a handful of runs, splits /
# merges a few runs so the proofErr boundaries land correctly, and
# removes three stray empty/page-break paragraphs.

$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function Get-ParaXmlOpen($range) {
    # Returns the opening <w:p ...> tag (with pPr) for the paragraph that
    # contains $range, captured via WordOpenXML so we can reuse rPr/pPr
    # untouched while only changing the run-level content.
    return $range.Paragraphs(1).Range
}

# ---------------------------------------------------------------------
# 1) "NORMAS Y PROCEDIMIENTOS A APLICAR" -> wrap the whole run in
#    proofErr gramStart/gramEnd
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("NORMAS Y PROCEDIMIENTOS A APLICAR", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p = Get-ParaXmlOpen $r
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="D9D9D9" w:themeFill="background1" w:themeFillShade="D9"/><w:spacing w:after="0" w:afterAutospacing="0"/><w:rPr><w:b/></w:rPr></w:pPr>' + `
  '<w:proofErr w:type="gramStart"/>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t>NORMAS Y PROCEDIMIENTOS A APLICAR</w:t></w:r>' + `
  '<w:proofErr w:type="gramEnd"/>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> (FAE)</w:t></w:r>' + `
  '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p.InsertXML($xml)

Write-Output "step1 done"

# ---------------------------------------------------------------------
# 2) "Asegurar que los entregables..." -> split into "Asegurar" (wrapped
#    in proofErr gramStart/gramEnd) + " que los entregables..."
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Asegurar que los entregables sean aprobados dentro de los plazos establecidos y cumplan con los estándares de calidad esperados.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p = $r.Paragraphs(1).Range
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="HelveticaNeueLT Std Med"/><w:iCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr>' + `
  '<w:proofErr w:type="gramStart"/>' + `
  '<w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="HelveticaNeueLT Std Med"/><w:iCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Asegurar</w:t></w:r>' + `
  '<w:proofErr w:type="gramEnd"/>' + `
  '<w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="HelveticaNeueLT Std Med"/><w:iCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> que los entregables sean aprobados dentro de los plazos establecidos y cumplan con los estándares de calidad esperados.</w:t></w:r>' + `
  '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p.InsertXML($xml)

Write-Output "step2 done"

# ---------------------------------------------------------------------
# 3) "EX.Y" -> wrap the whole run in proofErr gramStart/gramEnd
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("EX.Y", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p = $r.Paragraphs(1).Range
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:spacing w:afterAutospacing="0"/><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr>' + `
  '<w:proofErr w:type="gramStart"/>' + `
  '<w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>EX.Y</w:t></w:r>' + `
  '<w:proofErr w:type="gramEnd"/>' + `
  '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p.InsertXML($xml)

Write-Output "step3 done"

# ---------------------------------------------------------------------
# 4) Merge "Todas estas funcionalidades deberán" + " de estar
#    cumplimentadas..." into a single run (first paragraph, ends
#    with "...establecidos" and no trailing period)
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Todas estas funcionalidades deberán de estar cumplimentadas con sus correspondientes pruebas para probar que se llegan a los estándares de calidad establecidos", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p = $r.Paragraphs(1).Range
$rPr = '<w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="HelveticaNeueLT Std Med"/><w:iCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr>'
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="HelveticaNeueLT Std Med"/><w:iCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr>' + `
  "<w:r>$rPr" + '<w:t xml:space="preserve">En este entregable, se recogerán las funcionalidades relacionadas con </w:t></w:r>' + `
  "<w:r>$rPr" + '<w:t>la gestión de los cursos por parte del administrador</w:t></w:r>' + `
  "<w:r>$rPr" + '<w:t xml:space="preserve">, la visualización de los mismo por el cliente y su posterior reserva. </w:t></w:r>' + `
  "<w:r>$rPr" + '<w:t>Todas estas funcionalidades deberán de estar cumplimentadas con sus correspondientes pruebas para probar que se llegan a los estándares de calidad establecidos</w:t></w:r>' + `
  '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p.InsertXML($xml)

Write-Output "step4 done"

# ---------------------------------------------------------------------
# 5) Merge "la compra realizada con documentación enviada al usuario. "
#    + "Todas estas funcionalidades deberán de estar cumplimentadas..."
#    + "." into a single run
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("la compra realizada con documentación enviada al usuario.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p = $r.Paragraphs(1).Range
$rPr = '<w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="HelveticaNeueLT Std Med"/><w:iCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr>'
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="HelveticaNeueLT Std Med"/><w:iCs/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr>' + `
  "<w:r>$rPr" + '<w:t xml:space="preserve">En este entregable, se recogerán las funcionalidades relacionadas con </w:t></w:r>' + `
  "<w:r>$rPr" + '<w:t>la compra</w:t></w:r>' + `
  "<w:r>$rPr" + '<w:t xml:space="preserve"> de los cursos posterior a su reserva, incluyendo la pasarela de pago y la visualización de </w:t></w:r>' + `
  "<w:r>$rPr" + '<w:t>la compra realizada con documentación enviada al usuario. Todas estas funcionalidades deberán de estar cumplimentadas con sus correspondientes pruebas para probar que se llegan a los estándares de calidad establecidos.</w:t></w:r>' + `
  '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p.InsertXML($xml)

Write-Output "step5 done"

# ---------------------------------------------------------------------
# 6) "Actividad a incluir en el Diccionario de la EDT del paquete de
#    trabajo." -> split into "Actividad a incluir" (wrapped in proofErr
#    gramStart/gramEnd) + " en el Diccionario de la EDT del paquete de
#    trabajo."
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Actividad a incluir en el Diccionario de la EDT del paquete de trabajo.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p = $r.Paragraphs(1).Range
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:spacing w:afterAutospacing="0"/><w:jc w:val="both"/></w:pPr>' + `
  '<w:proofErr w:type="gramStart"/>' + `
  '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:i/><w:iCs/><w:color w:val="0070C0"/><w:sz w:val="20"/><w:lang w:val="es-ES"/></w:rPr><w:t>Actividad a incluir</w:t></w:r>' + `
  '<w:proofErr w:type="gramEnd"/>' + `
  '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:i/><w:iCs/><w:color w:val="0070C0"/><w:sz w:val="20"/><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> en el Diccionario de la EDT del paquete de trabajo.</w:t></w:r>' + `
  '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p.InsertXML($xml)

Write-Output "step6 done"

# ---------------------------------------------------------------------
# 7) "HITOS A INCLUIR" -> wrap the whole run in proofErr gramStart/gramEnd
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("HITOS A INCLUIR", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p = $r.Paragraphs(1).Range
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="HelveticaNeueLT Std Med"/><w:b/><w:bCs/></w:rPr></w:pPr>' + `
  '<w:proofErr w:type="gramStart"/>' + `
  '<w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="HelveticaNeueLT Std Med"/><w:b/><w:bCs/></w:rPr><w:t>HITOS A INCLUIR</w:t></w:r>' + `
  '<w:proofErr w:type="gramEnd"/>' + `
  '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p.InsertXML($xml)

Write-Output "step7 done"

# ---------------------------------------------------------------------
# 8) "REUNIONES A REALIZAR" -> wrap the whole run in proofErr
#    gramStart/gramEnd
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("REUNIONES A REALIZAR", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p = $r.Paragraphs(1).Range
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="HelveticaNeueLT Std Med"/><w:b/><w:bCs/></w:rPr></w:pPr>' + `
  '<w:proofErr w:type="gramStart"/>' + `
  '<w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="HelveticaNeueLT Std Med"/><w:b/><w:bCs/></w:rPr><w:t>REUNIONES A REALIZAR</w:t></w:r>' + `
  '<w:proofErr w:type="gramEnd"/>' + `
  '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p.InsertXML($xml)

Write-Output "step8 done"

# ---------------------------------------------------------------------
# 9) "INFORMES A ELABORAR" -> wrap the whole run in proofErr
#    gramStart/gramEnd
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("INFORMES A ELABORAR", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p = $r.Paragraphs(1).Range
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="HelveticaNeueLT Std Med"/><w:b/><w:bCs/></w:rPr></w:pPr>' + `
  '<w:proofErr w:type="gramStart"/>' + `
  '<w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="HelveticaNeueLT Std Med"/><w:b/><w:bCs/></w:rPr><w:t>INFORMES A ELABORAR</w:t></w:r>' + `
  '<w:proofErr w:type="gramEnd"/>' + `
  '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p.InsertXML($xml)

Write-Output "step9 done"
